$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.48
$ws.Range("B3").Value = 0.48
$ws.Range("B4").Value = 0.48
$ws.Range("C4").Value = 3
